# Update the "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting a re-scrape of the source site (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 8838
$ws1.Range("F10").Value = 84
$ws1.Range("F17").Value = 388
$ws1.Range("F18").Value = 11235
$ws1.Range("F25").Value = 146
$ws1.Range("F38").Value = 4103
$ws1.Range("F42").Value = 1261

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 12
$ws2.Range("F9").Value = 53

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 53
$ws4.Range("F11").Value = 8838
$ws4.Range("F13").Value = 84
$ws4.Range("F20").Value = 388
$ws4.Range("F21").Value = 11235
$ws4.Range("F24").Value = 146
$ws4.Range("F41").Value = 1261
